$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three menu-label cells at the top (A2-A4).
# Assigned in this order so the shared-string table appends new strings
# in the same order they appear in the target workbook.
$ws.Range("A2").Value = "🌐 Каталог TripTricks©"
$ws.Range("A3").Value = "💬 Задать вопрос"

# New survey-related rows (34-36), previously blank.
$ws.Range("A34").Value = "У вас закончились бесплатные ТрипТрики, можете пройти опрос для получения еще 5 ТрипТриков или пройти регистрацию по номеру телефона в личном кабинете и получить еще 5 ТрипТриков"
$ws.Rows.Item(34).RowHeight = 75

$ws.Range("A35").Value = "Пройти опрос"

$ws.Range("A36").Value = "Поздравляю, вы прошли опрос и получили свои 5 Трип Триков!"
$ws.Rows.Item(36).RowHeight = 30

# Personal-cabinet label, assigned last so it becomes the final new shared string.
$ws.Range("A4").Value = "🛫 Личный кабинет"

# Row 100 keeps its text but its height shrinks.
$ws.Rows.Item(100).RowHeight = 105

# Move the active selection from B1 to A4, as in the target workbook.
$ws.Range("A4").Select() | Out-Null
